$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2" = 0.8775636666666666
    "H2" = 2.632691
    "I2" = 0.1887436506618166
    "J2" = 0.2083714858314108
    "M2" = 63.46725166666666
    "N2" = 190.401755
    "O2" = 0.2354497988808272
    "P2" = 0.2397164477183668
    "Q2" = 55.69655408585611
    "R2" = 501.2689867727049
    "S2" = 0.04443965458835783
    "T2" = 0.04995007238930378
    "G3" = 0.8775636666666666
    "H3" = 2.632691
    "I3" = 0.1887436506618166
    "J3" = 0.2083714858314108
    "O3" = 0.1779985000094065
    "P3" = 0.1812240584798697
    "Q3" = 42.10622871669167
    "R3" = 378.956058450225
    "S3" = 0.03359608670410277
    "T3" = 0.03776192633384892
    "G4" = 0.8775636666666666
    "H4" = 2.632691
    "I4" = 0.1887436506618166
    "J4" = 0.2083714858314108
    "M4" = 64.53809233333334
    "N4" = 193.614277
    "O4" = 0.2394223865221556
    "P4" = 0.243761023683841
    "Q4" = 56.63628494771189
    "R4" = 509.726564529407
    "S4" = 0.04518945528235617
    "T4" = 0.05079284669278766
    "G5" = 0.8775636666666666
    "H5" = 2.632691
    "I5" = 0.1887436506618166
    "J5" = 0.2083714858314108
    "M5" = 14.3933435
    "N5" = 28.786687
    "O5" = 0.0533961963580272
    "P5" = 0.03624253541791403
    "Q5" = 12.63107529745283
    "R5" = 75.786451784717
    "S5" = 0.01007819303206925
    "T5" = 0.007551910955328276
    "G6" = 0.8775636666666666
    "H6" = 2.632691
    "I6" = 0.1887436506618166
    "J6" = 0.2083714858314108
    "M6" = 79.17795566666666
    "N6" = 237.533867
    "O6" = 0.2937331182295834
    "P6" = 0.2990559347000084
    "Q6" = 69.48369709401076
    "R6" = 625.353273846097
    "S6" = 0.05544026105493055
    "T6" = 0.0623147294601421
    "G7" = 2.458038666666667
    "H7" = 7.374116000000001
    "I7" = 0.5286672739959656
    "J7" = 0.5836444564186148
    "M7" = 63.46725166666666
    "N7" = 190.401755
    "O7" = 0.2354497988808272
    "P7" = 0.2397164477183668
    "Q7" = 156.0049586637311
    "R7" = 1404.04462797358
    "S7" = 0.1244746033372253
    "T7" = 0.1399091758231875
    "G8" = 2.458038666666667
    "H8" = 7.374116000000001
    "I8" = 0.5286672739959656
    "J8" = 0.5836444564186148
    "O8" = 0.1779985000094065
    "P8" = 0.1812240584798697
    "S8" = 0.09410198177534376
    "T8" = 0.1057704171014588
    "G9" = 2.458038666666667
    "H9" = 7.374116000000001
    "I9" = 0.5286672739959656
    "J9" = 0.5836444564186148
    "M9" = 64.53809233333334
    "N9" = 193.614277
    "O9" = 0.2394223865221556
    "P9" = 0.243761023683841
    "Q9" = 158.6371264282369
    "R9" = 1427.734137854132
    "S9" = 0.1265747804162764
    "T9" = 0.1422697701640005
    "G10" = 2.458038666666667
    "H10" = 7.374116000000001
    "I10" = 0.5286672739959656
    "J10" = 0.5836444564186148
    "M10" = 14.3933435
    "N10" = 28.786687
    "O10" = 0.0533961963580272
    "P10" = 0.03624253541791403
    "Q10" = 35.37939486561534
    "R10" = 212.276369193692
    "S10" = 0.02822882157035154
    "T10" = 0.02115275488322083
    "G11" = 2.458038666666667
    "H11" = 7.374116000000001
    "I11" = 0.5286672739959656
    "J11" = 0.5836444564186148
    "M11" = 79.17795566666666
    "N11" = 237.533867
    "O11" = 0.2937331182295834
    "P11" = 0.2990559347000084
    "Q11" = 194.6224765762858
    "R11" = 1751.602289186572
    "S11" = 0.1552870868967685
    "T11" = 0.1745423384467472
    "G12" = 1.313898
    "H12" = 2.627796
    "I12" = 0.2825890753422177
    "J12" = 0.2079840577499744
    "M12" = 63.46725166666666
    "N12" = 190.401755
    "O12" = 0.2354497988808272
    "P12" = 0.2397164477183668
    "Q12" = 83.38949503033
    "R12" = 500.33697018198
    "S12" = 0.0665355409552441
    "T12" = 0.04985719950587552
    "G13" = 1.313898
    "H13" = 2.627796
    "I13" = 0.2825890753422177
    "J13" = 0.2079840577499744
    "O13" = 0.1779985000094065
    "P13" = 0.1812240584798697
    "Q13" = 63.04191000585001
    "R13" = 378.2514600351
    "S13" = 0.0503004315299599
    "T13" = 0.03769171504456196
    "G14" = 1.313898
    "H14" = 2.627796
    "I14" = 0.2825890753422177
    "J14" = 0.2079840577499744
    "M14" = 64.53809233333334
    "N14" = 193.614277
    "O14" = 0.2394223865221556
    "P14" = 0.243761023683841
    "Q14" = 84.79647044058201
    "R14" = 508.778822643492
    "S14" = 0.06765815082352301
    "T14" = 0.05069840682705287
    "G15" = 1.313898
    "H15" = 2.627796
    "I15" = 0.2825890753422177
    "J15" = 0.2079840577499744
    "M15" = 14.3933435
    "N15" = 28.786687
    "O15" = 0.0533961963580272
    "P15" = 0.03624253541791403
    "Q15" = 18.911385237963
    "R15" = 75.645540951852
    "S15" = 0.0150891817556064
    "T15" = 0.007537869579364925
    "G16" = 1.313898
    "H16" = 2.627796
    "I16" = 0.2825890753422177
    "J16" = 0.2079840577499744
    "M16" = 79.17795566666666
    "N16" = 237.533867
    "O16" = 0.2937331182295834
    "P16" = 0.2990559347000084
    "Q16" = 104.031757594522
    "R16" = 624.1905455671319
    "S16" = 0.08300577027788428
    "T16" = 0.06219886679311912
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
